$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: subject 110 -> add first name, scan date, last name (needs rename later)
$ws.Range("B11").Value = "omri"

# "4.7.23" looks like a date to Excel's auto-detection; enter it as a formula
# that evaluates to the text, then paste the computed value back over itself
# (values only) so it is stored as a plain text/shared string, not a date
# serial number, and without picking up a new date number-format style.
$ws.Range("E11").Formula = '="4.7.23"'
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)

$ws.Range("C11").Value = "punaro"

# Row 12: subject 111 -> add first name, last name, scan date (needs rename later)
$ws.Range("B12").Value = "nataliya"
$ws.Range("C12").Value = "lukashov"

$ws.Range("E12").Formula = '="4.7.23"'
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)

# Update selection to match final cursor position
$ws.Range("E14").Select()
